# =====================================================================
# Applies the "APIInput.xlsx" update:
#   - refresh INPUT_SHEET test data (2 BANs instead of 4, new key words)
#   - drop the stray vertical "top" alignment that used to be applied
#     to the INPUT_SHEET BAN/API-key columns
#   - repoint TestURLinfo at the st2 API endpoints (keeping the display
#     text that was previously rendered for the two hyperlinked cells)
#   - duplicate the refreshed TestURLinfo tab as TestURLinfoST2
#   - hide subscription / Sheet1 / Sheet3 / Sheet2 tabs
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) INPUT_SHEET - new BAN / API key word pairs, trimmed to 2 rows
# ---------------------------------------------------------------
$inputSheet = $wb.Worksheets.Item("INPUT_SHEET")

$inputSheet.Range("A2").Value = 100040668
$inputSheet.Range("B2").Value = "TWREG-50571"
$inputSheet.Range("A3").Value = 100176138
$inputSheet.Range("B3").Value = "TWREG-50580"

# rows 4:6 (extra BANs + the "aslAccount,pastDue,billAmount" demo row)
# are no longer part of the sample data
$inputSheet.Rows("4:6").Delete()

# drop the vertical "top" alignment that was carried by this block -
# only horizontal "left" remains
$inputSheet.Range("A2:B3").VerticalAlignment = -4107

$inputSheet.Range("C9").Select()

# ---------------------------------------------------------------
# 2) TestURLinfo - point every API at the st2 test environment
# ---------------------------------------------------------------
$urlSheet = $wb.Worksheets.Item("TestURLinfo")

$urlSheet.Range("B2").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/process/sub/v1/accounts/$BAN/subscriptions'
$urlSheet.Range("B3").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/process/accounts/v1/accounts/$BAN/financial-status?realTimeInd=false&isGuestPay=false'
$urlSheet.Range("B6").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/process/sub/v1/accounts/$BAN/subscriptions/$SUBSCRIBER/contract'
$urlSheet.Range("B7").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/process/sub/v1/accounts/$BAN/subscriptions/$SUBSCRIBER/upgrade-eligibility?checkEarlyUpgrade=true&fetchAdditionalInfoKey=JUMPUPGRADE'
$urlSheet.Range("B8").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/process/v1/sub/accounts/$BAN/subscriptions/$SUBSCRIBER/current-services'
$urlSheet.Range("B9").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/process/eligibility/v1/flows?type=CHANGE_SERVICES&subscriberId=$SUBSCRIBER&accountId=$BAN&accountSubscriberIds=$SUBSCRIBER&role=ACCOUNT_OWNER'
$urlSheet.Range("B10").Value = 'https://st2-apiservices-sen.test.sprint.com:8442/api/digital/mac/v1/accounts/$BAN/adjustment-codes?macInd=autopay'

# B4 (future-payments) / B5 (payment-methods) are hyperlinked; the
# link itself keeps pointing at the st1 address, but since the cell
# text is moving to the st2 URL we must pin the previously-shown
# (st1) text as an explicit "display" override before rewriting the
# cell value
foreach ($hl in $urlSheet.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$5') {
        $hl.TextToDisplay = 'https://st1-apiservices-web.test.sprint.com:7441/api/process/pay/v1/accounts/$BAN/payment-methods?isUsgBan=false&isCheckAutoPayDiscount=true'
    }
    if ($addr -eq '$B$4') {
        $hl.TextToDisplay = 'https://st1-apiservices-web.test.sprint.com:7441/api/process/pay/v1/accounts/$BAN/future-payments?realTimeInd=false&isUsgBan=false'
    }
}

# now move the two hyperlinked cells over to the st2 URLs (TextToDisplay
# above also rewrote the cell text to the st1 string, so this is what
# actually lands in the cell)
$urlSheet.Range("B5").Value = 'https://st2-apiservices-web.test.sprint.com:7441/api/process/pay/v1/accounts/$BAN/payment-methods?isUsgBan=false&isCheckAutoPayDiscount=true'
$urlSheet.Range("B4").Value = 'https://st2-apiservices-web.test.sprint.com:7441/api/process/pay/v1/accounts/$BAN/future-payments?realTimeInd=false&isUsgBan=false'

$urlSheet.Range("A1:C10").Select()

# ---------------------------------------------------------------
# 3) Duplicate the refreshed TestURLinfo tab as TestURLinfoST2
# ---------------------------------------------------------------
$urlSheet.Copy([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "TestURLinfoST2"

# ---------------------------------------------------------------
# 4) Hide the helper / scratch tabs
# ---------------------------------------------------------------
$wb.Worksheets.Item("subscription").Visible = $false
$wb.Worksheets.Item("Sheet1").Visible = $false
$wb.Worksheets.Item("Sheet3").Visible = $false
$wb.Worksheets.Item("Sheet2").Visible = $false

$inputSheet.Select()
